$wb = $excel.ActiveWorkbook

# --- Update sigma_010 (sheet2) ---
$ws2 = $wb.Worksheets.Item("sigma_010")
$sigma010 = @(
    @(27.47287206360598, 29.39034561146602),
    @(27.48607003812207, 29.40452408891105),
    @(27.46750977065397, 29.39144744359502),
    @(27.47534051550214, 29.4204382957688),
    @(27.50380778484597, 29.38851878481415),
    @(27.49199222875973, 29.38716450881176),
    @(27.48167112512586, 29.4073266152364),
    @(27.49890340264041, 29.41207740506301),
    @(27.47016823339331, 29.36810542893682),
    @(27.4861094016698, 29.38025091686207),
    @(27.48344445643193, 29.39501990994651)
)
for ($i = 0; $i -lt $sigma010.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 2).Value = $sigma010[$i][0]
    $ws2.Cells.Item($row, 3).Value = $sigma010[$i][1]
}

# --- Update sigma_025 (sheet3) ---
$ws3 = $wb.Worksheets.Item("sigma_025")
$sigma025 = @(
    @(19.67075904166656, 25.69726061274941),
    @(19.6396086501729, 25.67420458294321),
    @(19.66471350038858, 25.72600560990867),
    @(19.67663969080323, 25.70926624709069),
    @(19.67388960355402, 25.70567620366846),
    @(19.64877050770563, 25.72126960805845),
    @(19.66652435848939, 25.71924214518252),
    @(19.68022021917821, 25.70413228440161),
    @(19.66858165866899, 25.72591061658974),
    @(19.65677417310814, 25.71339148473685),
    @(19.66464814037356, 25.70963593953296)
)
for ($i = 0; $i -lt $sigma025.Length; $i++) {
    $row = $i + 2
    $ws3.Cells.Item($row, 2).Value = $sigma025[$i][0]
    $ws3.Cells.Item($row, 3).Value = $sigma025[$i][1]
}

# --- Add new sheet sigma_050 (sheet4) at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "sigma_050"

$ws4.Cells.Item(1, 1).Value = "Rows"
$ws4.Cells.Item(1, 2).Value = "Noisy"
$ws4.Cells.Item(1, 3).Value = "NLM-LBP"

$sigma050 = @(
    @(0, 14.57963962169821, 21.0298455206133),
    @(1, 14.59158282359429, 20.9859510076853),
    @(2, 14.57148213461347, 21.03804352971977),
    @(3, 14.57091108084198, 21.02805087526399),
    @(4, 14.59786508174927, 21.02042117731074),
    @(5, 14.57319846358899, 21.08432407723135),
    @(6, 14.57499858504113, 21.046737724533),
    @(7, 14.57689288139636, 21.02508335977205),
    @(8, 14.59660533161504, 20.98853973205244),
    @(9, 14.56025688804831, 20.98658781864879)
)
for ($i = 0; $i -lt $sigma050.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 1).Value = $sigma050[$i][0]
    $ws4.Cells.Item($row, 2).Value = $sigma050[$i][1]
    $ws4.Cells.Item($row, 3).Value = $sigma050[$i][2]
}

$ws4.Cells.Item(12, 1).Value = "Média"
$ws4.Cells.Item(12, 2).Value = 14.57934328921871
$ws4.Cells.Item(12, 3).Value = 21.02335848228307

# Make sure active sheet ends up on first sheet like original
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
